# step_2 updated and reran
# Append a new "DiSCoVER: top drugs (cerebellar stem cell control)" slide
# (title textbox + results table) at the end of the deck. The new slide's
# content is identical to the existing slide 14, so we duplicate it and
# move the duplicate to the end of the deck — this reproduces the exact
# shape/table XML (including the default <p:clrMapOvr>) without relying
# on manually rebuilding every table cell via COM calls.

$p = $ppt.ActivePresentation

$source = $p.Slides.Item(14)

$dup = $source.Duplicate()
$newSlide = $dup.Item(1)
$newSlide.MoveTo($p.Slides.Count)

Write-Output "Slides: $($p.Slides.Count)"
